# Actualizar ligas y agregar Liga Argentina actualizada
# Append the newest Liga MX 2025 matchday rows (Jornada 7) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2025-07-19", "Necaxa",            "Club Queretaro",      3, 1, 1379491, 3,  2, 2, 0, 0, 0, 0, 0, 3, 1, "73%", "27%", "L"),
    @("2025-07-19", "Atletico San Luis",  "Monterrey",           0, 1, 1379492, 5,  3, 2, 2, 0, 0, 0, 0, 0, 1, "43%", "57%", "V"),
    @("2025-07-19", "Mazatlán",           "Puebla",              2, 1, 1379493, 3,  4, 3, 3, 0, 0, 0, 0, 2, 1, "42%", "58%", "L"),
    @("2025-07-20", "Tigres UANL",        "FC Juarez",           1, 0, 1379495, 5,  2, 3, 3, 0, 1, 0, 0, 1, 0, "64%", "36%", "L"),
    @("2025-07-20", "Leon",               "Guadalajara Chivas",  1, 0, 1379494, 5,  4, 1, 0, 0, 0, 0, 0, 1, 0, "33%", "67%", "L"),
    @("2025-07-20", "Atlas",              "Cruz Azul",           3, 3, 1379496, 4, 11, 2, 2, 0, 0, 0, 0, 3, 3, "28%", "72%", "E"),
    @("2025-07-20", "U.N.A.M. - Pumas",   "Pachuca",              2, 3, 1379497, 4,  1, 3, 2, 0, 0, 0, 0, 2, 3, "63%", "37%", "V")
)

$startRow = 12
$textColumns = @(1, 17, 18)   # A = Fecha, Q = Posesion Local (%), R = Posesion Visita (%) -- stored as literal text like "73%"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $col = $c + 1
        $cell = $ws.Cells.Item($r, $col)
        if ($textColumns -contains $col) {
            # Force text storage so values like "2025-07-19" / "73%" aren't
            # auto-converted into a date serial / numeric percentage.
            $cell.NumberFormat = "@"
            $cell.Value = $rowData[$c]
            $cell.ClearFormats()
        } else {
            $cell.Value = $rowData[$c]
        }
    }
}
